$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.4291351437568665
$ws.Range("B1").Value = 0.7030469179153442
$ws.Range("C1").Value = 2.199307680130005
$ws.Range("D1").Value = 4.683483123779297
$ws.Range("E1").Value = 2.162980556488037
